$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.939.34'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.87%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.119.67'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.91%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.92'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.624'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('B7').Value = 'Solana'
$ws.Range('C7').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.05'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.392'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0780'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.03%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.428.23'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.51'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.36'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.784'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.22'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.098.43'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.860.68'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.24'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.49'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0824'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.69'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.40%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.41'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.21'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.137'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +8.93%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.98'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.43'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.58'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.78%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.63'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.62'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0625'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.57'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.46'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +6.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.83'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.65%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.41'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.56%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.58%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.01'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.465.42'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0214'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.06'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.10%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.12'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -11.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.61'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.05'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.26'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.314.20'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.77%  '
